$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny float precision difference on existing row 13 date value
$ws.Range("A13").Value = 45813.39355503472

# Add new row 14 with updated price data
$ws.Range("A14").Value = 45814.39350229006
$ws.Range("B14").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C14").Value = "1Kg"
$ws.Range("D14").Value = "15,41€"

# Copy style of A13 (date column style) onto A14
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
